# Add pattern/annotation columns (E, G, H) to the "prompts-iterations"
# sheet, and make the corresponding style/bold adjustments on
# "mod_exp_summary" so that the (now-reused) style previously applied
# there keeps rendering bold while the new pattern cells render in the
# regular (non-bold) font.

$wb = $excel.ActiveWorkbook

$wsPrompts = $wb.Worksheets.Item("prompts-iterations")
$wsSummary = $wb.Worksheets.Item("mod_exp_summary")

# --- New "pattern" columns on prompts-iterations ---------------------
# Row 2 (entry "ftexbqmj" / prompt1)
$wsPrompts.Range("E2").Value = "parallelize"
$wsPrompts.Range("G2").Value = "parallel"
$wsPrompts.Range("H2").Value = "A9"

# Row 3 (prompt2)
$wsPrompts.Range("E3").Value = "loop"
$wsPrompts.Range("G3").Value = "loop"
$wsPrompts.Range("H3").Value = "A8"

# Row 4 (prompt3)
$wsPrompts.Range("E4").Value = "condition"
$wsPrompts.Range("G4").Value = "condition"
$wsPrompts.Range("H4").Value = "A13"

# Row 5 (new row, previously empty)
$wsPrompts.Range("G5").Value = "replace"
$wsPrompts.Range("H5").Value = "A4"

# Row 6 (entry "qhslybdt" / prompt1)
$wsPrompts.Range("E6").Value = "loop"

# Row 7 (prompt2)
$wsPrompts.Range("E7").Value = "replace"

# Row 8 (prompt3)
$wsPrompts.Range("E8").Value = "loop"

# Row 10 (entry "qzicwyto" / prompt1)
$wsPrompts.Range("E10").Value = "data objects"

# Row 12 (entry "rqvtdhws" / prompt1)
$wsPrompts.Range("E12").Value = "parallelize"

# Row 13 (prompt2)
$wsPrompts.Range("E13").Value = "parallelize"

# Row 15 (entry "rxsvqgua" / prompt1)
$wsPrompts.Range("E15").Value = "parallelize"

# Row 17 (entry "cjaktieq" / prompt1)
$wsPrompts.Range("E17").Value = "parallelize"

# Row 18 (prompt2)
$wsPrompts.Range("E18").Value = "loop"

# Row 19 (prompt3)
$wsPrompts.Range("E19").Value = "loop"

# New pattern cells use the regular (non-bold) font.
$wsPrompts.Range("E2:E4,E6:E8,E10,E12:E13,E15,E17:E19").Font.Bold = $false
$wsPrompts.Range("G2:G5").Font.Bold = $false
$wsPrompts.Range("H2:H5").Font.Bold = $false

# --- mod_exp_summary keeps its bold look on N14:O15 -------------------
$wsSummary.Range("N14:O15").Font.Bold = $true

# --- selection / active sheet -----------------------------------------
$wsPrompts.Activate()
$wsPrompts.Range("G16").Select()
